# Auto-generated: route the 5 "Hub" shapes into a single group ("Group 51")
# matching the target OOXML (grpSp off/ext/chOff/chExt).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Step 1: temporarily move the 5 shapes to their rotation-aware visual bounding
# boxes. iron_native computes a new groups child coordinate space (chOff/chExt) as the
# plain union of each members raw Left/Top/Width/Height (which for a rotated shape is
# its un-rotated off/ext, NOT the rotated visual bbox that real PowerPoint uses). Faking
# the visual bbox here - then restoring true values after Group() - reproduces the
# chOff/chExt real PowerPoint would have written.

# TextBox10 (shape #9): visual bbox off=(6094301,2413969) ext=(3377946,1200329)
$c9 = $s.Shapes.Item(9)
$c9.Left = 479.8662261968504
$c9.Top = 190.07630157480315
$c9.Width = 265.98
$c9.Height = 94.51409530708662

# LeftBracket12 (shape #10): visual bbox off=(6088150,3185157) ext=(3377943,429139)
$c10 = $s.Shapes.Item(10)
$c10.Left = 479.38188976377955
$c10.Top = 250.7997665433071
$c10.Width = 265.97976377952756
$c10.Height = 33.79047244094488

# LeftBracket14 (shape #11): visual bbox off=(5884736,2895992) ext=(3748161,971632)
$c11 = $s.Shapes.Item(11)
$c11.Left = 463.36503937007876
$c11.Top = 228.03086614173228
$c11.Width = 295.1307874015748
$c11.Height = 76.50645828346457

# LeftBracket16 (shape #12): visual bbox off=(5679482,3185158) ext=(4195276,971629)
$c12 = $s.Shapes.Item(12)
$c12.Left = 447.2033070866142
$c12.Top = 250.79984283464566
$c12.Width = 330.33669291338583
$c12.Height = 76.50622177165354

# LeftBracket17 (shape #13): visual bbox off=(5472204,3032756) ext=(4609834,1418211)
$c13 = $s.Shapes.Item(13)
$c13.Left = 430.8822174094488
$c13.Top = 238.7996902519685
$c13.Width = 362.97905511811024
$c13.Height = 111.67015748031496

# --- Step 2: group the (temporarily repositioned) shapes.
$range = $s.Shapes.Range(@(9,10,11,12,13))
$grp = $range.Group()

# --- Step 3: restore each child shape to its true original Left/Top/Width/Height.
# (Group-level chOff/chExt were already captured at Group() time and are unaffected
# by subsequent child-shape moves.)
# TextBox10: restore off=(6094301,2413969) ext=(3377946,1200329)
$gc9 = $grp.GroupItems.Item(1)
$gc9.Left = 479.8662261968504
$gc9.Top = 190.07630157480315
$gc9.Width = 265.98
$gc9.Height = 94.51409530708662

# LeftBracket12: restore off=(7562552,1710755) ext=(429139,3377943)
$gc10 = $grp.GroupItems.Item(2)
$gc10.Left = 595.4765354330709
$gc10.Top = 134.7051181102362
$gc10.Width = 33.79047244094488
$gc10.Height = 265.97976377952756

# LeftBracket14: restore off=(7273000,1507728) ext=(971632,3748161)
$gc11 = $grp.GroupItems.Item(3)
$gc11.Left = 572.6771653543307
$gc11.Top = 118.71874015748031
$gc11.Width = 76.50645828346457
$gc11.Height = 295.1307874015748

# LeftBracket16: restore off=(7291306,1573335) ext=(971629,4195276)
$gc12 = $grp.GroupItems.Item(4)
$gc12.Left = 574.1185826771654
$gc12.Top = 123.88464737007874
$gc12.Width = 76.50622177165354
$gc12.Height = 330.33669291338583

# LeftBracket17: restore off=(7068015,1436945) ext=(1418211,4609834)
$gc13 = $grp.GroupItems.Item(5)
$gc13.Left = 556.5366141732284
$gc13.Top = 113.14527559055118
$gc13.Width = 111.67015748031496
$gc13.Height = 362.97905511811024

# --- Step 4: set the groups final outer position/size on the slide.
$grp.Left = 356.41401672440946
$grp.Top = 145.89472440944883
$grp.Width = 437.44725037007873
$grp.Height = 204.57511811023622

# --- Step 5: name the new group to match the authored deck.
$grp.Name = "Group 51"
